{"js": "// Update the date line at the top of the document.\nconst dateResults = context.document.body.search(\"2024-03-15 Friday\", { matchCase: true });\ndateResults.load('items');\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-03-16 Saturday\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Update the division problems in the practice table, addressed by (row, column)\n// so that duplicate text values elsewhere in the table can never cause a mismatch.\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  { row: 0, col: 0, text: \"17\u00f77=\" },\n  { row: 0, col: 1, text: \"58\u00f77=\" },\n  { row: 0, col: 2, text: \"65\u00f73=\" },\n  { row: 0, col: 3, text: \"24\u00f73=\" },\n  { row: 0, col: 4, text: \"55\u00f79=\" },\n  { row: 4, col: 0, text: \"24\u00f75=\" },\n  { row: 4, col: 1, text: \"52\u00f74=\" },\n  { row: 4, col: 2, text: \"64\u00f77=\" },\n  { row: 4, col: 3, text: \"45\u00f77=\" },\n  { row: 4, col: 4, text: \"75\u00f76=\" },\n  { row: 8, col: 0, text: \"79\u00f75=\" },\n  { row: 8, col: 1, text: \"76\u00f72=\" },\n  { row: 8, col: 2, text: \"54\u00f75=\" },\n  { row: 8, col: 3, text: \"44\u00f76=\" },\n  { row: 8, col: 4, text: \"35\u00f79=\" },\n  { row: 12, col: 0, text: \"93\u00f74=\" },\n  { row: 12, col: 1, text: \"10\u00f72=\" },\n  { row: 12, col: 2, text: \"86\u00f78=\" },\n  { row: 12, col: 3, text: \"43\u00f75=\" },\n  { row: 12, col: 4, text: \"60\u00f76=\" },\n  { row: 16, col: 0, text: \"17\u00f75=\" },\n  { row: 16, col: 1, text: \"42\u00f76=\" },\n  { row: 16, col: 2, text: \"56\u00f78=\" },\n  { row: 16, col: 3, text: \"54\u00f73=\" },\n  { row: 16, col: 4, text: \"78\u00f74=\" },\n];\n\nfor (const update of cellUpdates) {\n  const cell = table.getCell(update.row, update.col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load('items');\n  await context.sync();\n  const range = paragraphs.items[0].getRange();\n  range.insertText(update.text, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the document.\n$dateFind = $d.Content.Find\n$dateFind.ClearFormatting()\n$dateFind.Replacement.ClearFormatting()\n$dateFind.Text = \"2024-03-15 Friday\"\n$dateFind.Replacement.Text = \"2024-03-16 Saturday\"\n$dateFind.Execute($null,$false,$false,$false,$false,$false,$true,1,$false,$null,2) | Out-Null\n\n# Update the division problems in the practice table, addressed by (row, column)\n# so that duplicate text values elsewhere in the table can never cause a mismatch.\n$tbl = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; Text = \"17\u00f77=\" }\n    @{ Row = 1; Col = 2; Text = \"58\u00f77=\" }\n    @{ Row = 1; Col = 3; Text = \"65\u00f73=\" }\n    @{ Row = 1; Col = 4; Text = \"24\u00f73=\" }\n    @{ Row = 1; Col = 5; Text = \"55\u00f79=\" }\n    @{ Row = 5; Col = 1; Text = \"24\u00f75=\" }\n    @{ Row = 5; Col = 2; Text = \"52\u00f74=\" }\n    @{ Row = 5; Col = 3; Text = \"64\u00f77=\" }\n    @{ Row = 5; Col = 4; Text = \"45\u00f77=\" }\n    @{ Row = 5; Col = 5; Text = \"75\u00f76=\" }\n    @{ Row = 9; Col = 1; Text = \"79\u00f75=\" }\n    @{ Row = 9; Col = 2; Text = \"76\u00f72=\" }\n    @{ Row = 9; Col = 3; Text = \"54\u00f75=\" }\n    @{ Row = 9; Col = 4; Text = \"44\u00f76=\" }\n    @{ Row = 9; Col = 5; Text = \"35\u00f79=\" }\n    @{ Row = 13; Col = 1; Text = \"93\u00f74=\" }\n    @{ Row = 13; Col = 2; Text = \"10\u00f72=\" }\n    @{ Row = 13; Col = 3; Text = \"86\u00f78=\" }\n    @{ Row = 13; Col = 4; Text = \"43\u00f75=\" }\n    @{ Row = 13; Col = 5; Text = \"60\u00f76=\" }\n    @{ Row = 17; Col = 1; Text = \"17\u00f75=\" }\n    @{ Row = 17; Col = 2; Text = \"42\u00f76=\" }\n    @{ Row = 17; Col = 3; Text = \"56\u00f78=\" }\n    @{ Row = 17; Col = 4; Text = \"54\u00f73=\" }\n    @{ Row = 17; Col = 5; Text = \"78\u00f74=\" }\n)\n\nforeach ($update in $cellUpdates) {\n    $cell = $tbl.Cell($update.Row, $update.Col)\n    $cell.Range.Text = $update.Text\n}\n\nWrite-Output \"done\"\n"}
